# Fix row 21 in Oaks Crisis
# Also corrects related downstream totals in rows 4, 8, 13, 14, 18, 28, 29

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J4").Value = 45
$ws.Range("Q4").Value = 301

$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 40
$ws.Range("Q8").Value = 216

$ws.Range("K13").Value = 26
$ws.Range("Q13").Value = 82

$ws.Range("J14").Value = 76
$ws.Range("K14").Value = 81
$ws.Range("Q14").Value = 356

$ws.Range("J18").Value = 15
$ws.Range("Q18").Value = 110

$ws.Range("J21").Value = 74
$ws.Range("K21").Value = 72
$ws.Range("Q21").Value = 269

$ws.Range("J28").Value = 11
$ws.Range("Q28").Value = 95

$ws.Range("J29").Value = 27
$ws.Range("Q29").Value = 72

$wb.Save()
